# Adds a "Request/sec" column (computed as existing-throughput * 60) to each
# of the V1-V4 report sheets, and updates the active-sheet/selection state
# to match the author's final view (V1 active, cell I6 selected there).

$wb = $excel.ActiveWorkbook

foreach ($name in @("V1", "V2", "V3", "V4")) {
    $ws = $wb.Worksheets.Item($name)

    # First summary block (row 2 header, rows 3/4 data)
    $ws.Range("L2").Value = "Request/sec"
    $ws.Range("L3").Formula = "=H3*60"
    $ws.Range("L4").Formula = "=H4*60"

    # Second "Aggregate" table (row 9 header, rows 10/11 data)
    $ws.Range("N9").Value = "Request/sec"
    $ws.Range("N10").Formula = "=K10*60"
    $ws.Range("N11").Formula = "=K11*60"
}

# Restore per-sheet selections (selecting a range also activates that
# sheet, so do the non-active sheets first and finish on V1 so it ends
# up the active tab, matching the target workbook view state).
$wb.Worksheets.Item("V2").Range("L3").Select()
$wb.Worksheets.Item("V3").Range("L4").Select()
$wb.Worksheets.Item("V4").Range("L2:L4").Select()
$wb.Worksheets.Item("V1").Range("I6").Select()
